# Revisions required for new jlcpcb submission
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header column G from "Side" to "Layer" (data underneath is unchanged)
$ws.Range("G1").Value = "Layer"

# 2. The header row no longer carries the bold "title" style - clear it back
#    to the sheet default formatting.
$ws.Rows.Item(1).ClearFormats()

# 3. U1 (row 35) now references a new custom footprint instead of the stock
#    LQFP-64 package.
$ws.Range("C35").Value = "my_STM32F405RGTx_2"

# 4. U1's placement rotation changed from 180 to -90 degrees.
$ws.Range("F35").Value = -90

# 5. Leave the active selection on G1, matching the saved view state.
$ws.Range("G1").Select()
